$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text fields (coin names, links) - assign directly.
$textCells = @(
    @{ Addr = 'B19'; Value = 'ShibaInu' },
    @{ Addr = 'C19'; Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib' },
    @{ Addr = 'B20'; Value = 'WrappedliquidstakedEther2.0' },
    @{ Addr = 'C20'; Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth' },
    @{ Addr = 'B51'; Value = 'Elrond' },
    @{ Addr = 'C51'; Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld' }
)

# Price / volume fields look numeric (e.g. "1.000", "30.530.22", "  +0.46%  ")
# but must stay plain text, matching the source data feeds formatting.
# Force the cell to Text format first so Excel does not reinterpret the
# string as a number/date and strip the literal formatting.
$numericLookingCells = @(
    @{ Addr = 'D2'; Value = '30.530.22' },
    @{ Addr = 'E2'; Value = '  +0.46%  ' },
    @{ Addr = 'D3'; Value = '1.916.26' },
    @{ Addr = 'E3'; Value = '  +0.01%  ' },
    @{ Addr = 'D4'; Value = '1.000' },
    @{ Addr = 'E4'; Value = '  -0.16%  ' },
    @{ Addr = 'D5'; Value = '244.37' },
    @{ Addr = 'E5'; Value = '  +1.31%  ' },
    @{ Addr = 'D6'; Value = '1.000' },
    @{ Addr = 'E6'; Value = '  -0.13%  ' },
    @{ Addr = 'D7'; Value = '0.4887' },
    @{ Addr = 'E7'; Value = '  +4.50%  ' },
    @{ Addr = 'D8'; Value = '0.2905' },
    @{ Addr = 'E8'; Value = '  +2.37%  ' },
    @{ Addr = 'D9'; Value = '0.06713' },
    @{ Addr = 'E9'; Value = '  -2.76%  ' },
    @{ Addr = 'D10'; Value = '106.88' },
    @{ Addr = 'E10'; Value = '  +1.60%  ' },
    @{ Addr = 'D11'; Value = '18.84' },
    @{ Addr = 'E11'; Value = '  +3.88%  ' },
    @{ Addr = 'D12'; Value = '1.923.47' },
    @{ Addr = 'E12'; Value = '  +0.48%  ' },
    @{ Addr = 'D13'; Value = '0.07621' },
    @{ Addr = 'E13'; Value = '  -0.34%  ' },
    @{ Addr = 'D14'; Value = '5.273' },
    @{ Addr = 'E14'; Value = '  +2.54%  ' },
    @{ Addr = 'D15'; Value = '0.6653' },
    @{ Addr = 'E15'; Value = '  +1.90%  ' },
    @{ Addr = 'D16'; Value = '273.22' },
    @{ Addr = 'E16'; Value = '  -3.81%  ' },
    @{ Addr = 'D17'; Value = '30.521.67' },
    @{ Addr = 'E17'; Value = '  +0.47%  ' },
    @{ Addr = 'D18'; Value = '1.000' },
    @{ Addr = 'E18'; Value = '  -0.03%  ' },
    @{ Addr = 'D19'; Value = '0.000007532' },
    @{ Addr = 'E19'; Value = '  -0.55%  ' },
    @{ Addr = 'D20'; Value = '2.168.30' },
    @{ Addr = 'E20'; Value = '  +0.19%  ' },
    @{ Addr = 'D21'; Value = '12.82' },
    @{ Addr = 'E21'; Value = '  -0.87%  ' },
    @{ Addr = 'D22'; Value = '5.492' },
    @{ Addr = 'E22'; Value = '  +5.57%  ' },
    @{ Addr = 'D23'; Value = '0.9997' },
    @{ Addr = 'E23'; Value = '  -0.34%  ' },
    @{ Addr = 'D24'; Value = '6.396' },
    @{ Addr = 'E24'; Value = '  +3.78%  ' },
    @{ Addr = 'D25'; Value = '9.406' },
    @{ Addr = 'E25'; Value = '  +2.09%  ' },
    @{ Addr = 'D26'; Value = '163.78' },
    @{ Addr = 'E26'; Value = '  -2.47%  ' },
    @{ Addr = 'D27'; Value = '20.03' },
    @{ Addr = 'E27'; Value = '  -5.25%  ' },
    @{ Addr = 'D28'; Value = '2.107' },
    @{ Addr = 'E28'; Value = '  +4.02%  ' },
    @{ Addr = 'D29'; Value = '0.1049' },
    @{ Addr = 'E29'; Value = '  -1.73%  ' },
    @{ Addr = 'D30'; Value = '1.402' },
    @{ Addr = 'E30'; Value = '  +2.29%  ' },
    @{ Addr = 'D31'; Value = '4.123' },
    @{ Addr = 'E31'; Value = '  +0.38%  ' },
    @{ Addr = 'D32'; Value = '4.048' },
    @{ Addr = 'E32'; Value = '  +2.62%  ' },
    @{ Addr = 'D33'; Value = '0.04995' },
    @{ Addr = 'E33'; Value = '  -1.41%  ' },
    @{ Addr = 'D34'; Value = '0.7264' },
    @{ Addr = 'E34'; Value = '  -1.09%  ' },
    @{ Addr = 'D35'; Value = '1.135' },
    @{ Addr = 'E35'; Value = '  -0.52%  ' },
    @{ Addr = 'E36'; Value = '  +0.10%  ' },
    @{ Addr = 'D37'; Value = '2.726' },
    @{ Addr = 'E37'; Value = '  -0.26%  ' },
    @{ Addr = 'D38'; Value = '0.02031' },
    @{ Addr = 'E38'; Value = '  +1.19%  ' },
    @{ Addr = 'D39'; Value = '2.672' },
    @{ Addr = 'E39'; Value = '  -0.24%  ' },
    @{ Addr = 'D40'; Value = '110.82' },
    @{ Addr = 'E40'; Value = '  +2.56%  ' },
    @{ Addr = 'D41'; Value = '2.017' },
    @{ Addr = 'E41'; Value = '  -1.19%  ' },
    @{ Addr = 'D42'; Value = '0.4424' },
    @{ Addr = 'E42'; Value = '  +5.88%  ' },
    @{ Addr = 'D43'; Value = '0.8657' },
    @{ Addr = 'E43'; Value = '  -0.52%  ' },
    @{ Addr = 'D44'; Value = '5.881' },
    @{ Addr = 'E44'; Value = '  +1.18%  ' },
    @{ Addr = 'D45'; Value = '1.000' },
    @{ Addr = 'E45'; Value = '  -0.09%  ' },
    @{ Addr = 'D46'; Value = '67.87' },
    @{ Addr = 'E46'; Value = '  +1.15%  ' },
    @{ Addr = 'D47'; Value = '7.299' },
    @{ Addr = 'E47'; Value = '  +2.33%  ' },
    @{ Addr = 'D48'; Value = '9.323' },
    @{ Addr = 'E48'; Value = '  +1.49%  ' },
    @{ Addr = 'D49'; Value = '0.1245' },
    @{ Addr = 'E49'; Value = '  +3.93%  ' },
    @{ Addr = 'D50'; Value = '47.70' },
    @{ Addr = 'E50'; Value = '  -8.23%  ' },
    @{ Addr = 'D51'; Value = '34.58' },
    @{ Addr = 'E51'; Value = '  +0.24%  ' }
)

foreach ($item in $textCells) {
    $ws.Range($item.Addr).Value = $item.Value
}

foreach ($item in $numericLookingCells) {
    $range = $ws.Range($item.Addr)
    $range.NumberFormat = "@"
    $range.Value = $item.Value
}